$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns I:L entirely (shrinks dimension to A1:H5)
$ws.Range("I1:L5").Delete()

# Row 1 header cells become raw numeric values (course-name parsing bug)
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 1234567891
$ws.Range("F1").Value = 3322111234
$ws.Range("G1").Value = 4453245321
$ws.Range("H1").Value = 5555555555

# Row 2 (John Doe)
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# Row 3 (Andrew Hartmann)
$ws.Range("F3").Value = "LATE (completed)"
$ws.Range("H3").Value = ""

# Row 4 (Nick -> Nicholas Fletcher)
$ws.Range("C4").Value = "Nicholas"
$ws.Range("G4").Value = "Completed"
$ws.Range("H4").Value = ""

# Row 5 (John -> J Cena)
$ws.Range("C5").Value = "J"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = "LATE (completed)"
